$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the new block of paragraphs after "Site_weather_Dep ..." and
#    before "Nutrient_data, csv file format for use"
# ---------------------------------------------------------------------------

$GREEN = 5287936  # RGB(0x00, 0xB0, 0x50) packed as BGR for Word's Font.Color

$items = @(
    @{Text="Process_NETCDF_Data_For_Use.R "; Level=3; Color=$true},
    @{Text="## This script process NETCDF data output from CMAQ bidi "; Level=4; Color=$false},
    @{Text="## Format:"; Level=4; Color=$false},
    @{Text="## Set Up lists of variables and File names"; Level=4; Color=$false},
    @{Text="## Function to pull them from NETCDF (Pull_Var)"; Level=4; Color=$false},
    @{Text="## Function to convert Var to data frame (Mat_to_Data_Frame)"; Level=4; Color=$false},
    @{Text="## Loop to apply Pull_Var and use same Loop to apply Mat_to_Data_Frame *applys to all variables "; Level=4; Color=$false},
    @{Text="# previously applied to Pull_Var because a list of variables was created in the Pull_Var loop"; Level=4; Color=$false},
    @{Text="## CSV FILES WILL BE FOUND IN WORKING DIRECTORY "; Level=4; Color=$false},
    @{Text="Process_Individual_CSV_Var_Data_Over_Watershed.R"; Level=3; Color=$true},
    @{Text="## This script takes the watersheds and aggragates the variables across each and then writes a CSV file for each variable that includes the"; Level=4; Color=$false},
    @{Text="# new number, date, and watershed id"; Level=4; Color=$false}
)

# Find the anchor paragraph ("Site_weather_Dep ...") - new paragraphs are
# inserted right after it, one at a time, so they stay in order.
$anchor = $d.Paragraphs(4)

$coloredParas = @()

foreach ($it in $items) {
    $anchor.Range.InsertParagraphAfter()
    $newPara = $anchor.Next()
    $newPara.Range.Text = $it.Text
    $newPara.Range.ListFormat.ListLevelNumber = $it.Level
    if ($it.Color) {
        $coloredParas += $newPara
    }
    $anchor = $newPara
}

foreach ($cp in $coloredParas) {
    $cp.Range.Font.Color = $GREEN
}

Write-Output "Inserted block OK"

# ---------------------------------------------------------------------------
# 2. Colour the "Process_CSV_Nutrient_Data.R " run green (only that run, the
#    rest of the paragraph is untouched).
# ---------------------------------------------------------------------------

$rng = $d.Content
$found = $rng.Find.Execute("Process_CSV_Nutrient_Data.R ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Font.Color = $GREEN
}
Write-Output "Colored Process_CSV_Nutrient_Data.R run: $found"

# ---------------------------------------------------------------------------
# 3. Colour the whole "Load_Dataframe", "GLS_With_Function" and "CV_Custom"
#    paragraphs green (paragraph mark + run).
# ---------------------------------------------------------------------------

$wholeParaTargets = @("Load_Dataframe", "GLS_With_Function", "CV_Custom")
foreach ($target in $wholeParaTargets) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs($i)
        if ($para.Range.Text.Trim() -eq $target) {
            $para.Range.Font.Color = $GREEN
            Write-Output "  colored paragraph $i ($target)"
            break
        }
    }
}
Write-Output "Colored whole paragraphs OK"

# ---------------------------------------------------------------------------
# 4. Rearrange the "Special note ..." text. The words "portion, if it is not
#    make sure" move from the start of the 2nd run to the end of the 1st run
#    (the _GoBack bookmark sits between the two runs and must stay put).
# ---------------------------------------------------------------------------

$r1 = $d.Content
$found2 = $r1.Find.Execute("Special note, if n is a factor of total sample size comment out the remainder", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Found run1 anchor: $found2 [$($r1.Text)] start=$($r1.Start) end=$($r1.End)"

$r1Start = $r1.Start
$r1Fixed = $d.Range($r1Start, $r1.End)
$r1Fixed.Text = "Special note, if n is a factor of total sample size comment out the remainder portion, if it is not make sure"
Write-Output "Run1 updated: [$($r1Fixed.Text)] start=$($r1Fixed.Start) end=$($r1Fixed.End)"

$moveLen = (" portion, if it is not make sure").Length
$trimStart = $r1Fixed.End
$r2trim = $d.Range($trimStart, $trimStart + $moveLen)
Write-Output "About to trim: [$($r2trim.Text)]"
$r2trim.Text = ""

Write-Output "Special note fix OK"
